# Corrected most names to the official names from website
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# District name corrections in column G
$ws.Range("G4").Value = "Davangere"
$ws.Range("G9").Value = "Davangere"
$ws.Range("G10").Value = "Davangere"
$ws.Range("G13").Value = "Davangere"
$ws.Range("G14").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G24").Value = "Ramanagara"
$ws.Range("G26").Value = "Bagalkot"
$ws.Range("G28").Value = "Vijayapura (Bijapur)"
$ws.Range("G29").Value = "Chikkamagaluru (Chikmagalur)"
$ws.Range("G30").Value = "Vijayapura (Bijapur)"
$ws.Range("G31").Value = "Kalaburagi (Gulbarga)"

# Remove now-empty stray cells in column F (clearing collapses them out of the sheet)
$ws.Range("F6").ClearContents() | Out-Null
$ws.Range("F17").ClearContents() | Out-Null
$ws.Range("F33").ClearContents() | Out-Null
